# Apply "normalized/denormalized layout" edit to Collection_SC sheet.
# For the EDULEVEL and EDUYRNUM collection groups, add a parallel
# "Denormalized" implementation_option row-set: column H becomes the
# *_DENORMALIZED collection_group_id, column I gets the new
# implementation_option value "Denormalized", column J gets the new
# scenario number 360, and column K (short_name) gets a "(Denormalized)"
# suffix appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Collection_SC")

# Fill column-by-column (mirrors a fill-down workflow: implementation_option
# first, then collection_group_id, then scenario, then short_name) so the
# shared-string table ends up in the same append order as the source edit.
$ws.Range("I2:I6").Value = "Denormalized"

$ws.Range("H2").Value = "EDULEVEL_DENORMALIZED"
$ws.Range("H3").Value = "EDULEVEL_DENORMALIZED"
$ws.Range("H4").Value = "EDUYRNUM_DENORMALIZED"
$ws.Range("H5").Value = "EDUYRNUM_DENORMALIZED"
$ws.Range("H6").Value = "EDUYRNUM_DENORMALIZED"

$ws.Range("J2:J6").Value = 360

$ws.Range("K2").Value = "Education Level (Denormalized)"
$ws.Range("K3").Value = "Education Level (Denormalized)"
$ws.Range("K4").Value = "Number of Years of Education (Denormalized)"
$ws.Range("K5").Value = "Number of Years of Education (Denormalized)"
$ws.Range("K6").Value = "Number of Years of Education (Denormalized)"

# Column width / layout touch-ups matching the refreshed sheet formatting.
$ws.Columns.Item(8).ColumnWidth = 15.83
$ws.Columns.Item(9).ColumnWidth = 17
$ws.Columns.Item(3).ColumnWidth = 22.66
$ws.Columns.Item(11).ColumnWidth = 40.16
$ws.Columns.Item(30).ColumnWidth = 19.33
$ws.Columns.Item(33).ColumnWidth = 15.33

# Match the final selection left in the sheet.
$ws.Range("I19").Select() | Out-Null
